$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 13
$ws.Range("AG2").Value = 12

# Row 3 updates
$ws.Range("I3").Value = 1.5
$ws.Range("L3").Value = 1.25
$ws.Range("N3").Value = 4.9
$ws.Range("R3").Value = 1.54
$ws.Range("S3").Value = 2.48
$ws.Range("T3").Value = 1.84
$ws.Range("U3").Value = 1.94
$ws.Range("V3").Value = 3
$ws.Range("X3").Value = 29
$ws.Range("Y3").Value = 12.5
$ws.Range("Z3").Value = 11
$ws.Range("AA3").Value = 15
$ws.Range("AB3").Value = 38
$ws.Range("AC3").Value = 15
$ws.Range("AD3").Value = 12.5
$ws.Range("AE3").Value = 18
$ws.Range("AG3").Value = 40
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 40
$ws.Range("AO3").Value = 6.8

# Row 4 updates
$ws.Range("N4").Value = 1.34
$ws.Range("P4").Value = 1.34
$ws.Range("S4").Value = 1.2
